$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 6) with the "best" model results.
$ws.Cells.Item(6, 1).Value = 44140.774085648147
$ws.Cells.Item(6, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(6, 2).Value = "lgbm3"
$ws.Cells.Item(6, 3).Value = "best"
$ws.Cells.Item(6, 4).Value = "best.trn.feather"
$ws.Cells.Item(6, 5).Value = "rmse"
$ws.Cells.Item(6, 6).Value = 1.92169486412754
$ws.Cells.Item(6, 7).Value = "lgbm3.best.sub.csv"
$ws.Cells.Item(6, 8).Value = 0.60338999999999998

$ws.Range("H6").Select()
